$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 4 with more meaningful values
$ws.Range("B4").Value = "Una definición más propia"
$ws.Range("C4").Value = "Una regla más propia"
$ws.Range("D4").Value = "Palabra_1, Palabra_2"

# Add new row 7: Calle
$ws.Range("A7").Value = "Calle"
$ws.Range("B7").Value = "iojfiwjfifj de la calle"
$ws.Range("C7").Value = "iojargiwjfpowfdpo"
$ws.Range("D7").Value = "street, otro"

# Add new row 8: Monto
$ws.Range("A8").Value = "Monto"
$ws.Range("B8").Value = "calculo que seiuhnfowijef"
$ws.Range("C8").Value = "varianza de cartera + varianza de clientes = bono"
$ws.Range("D8").Value = "comisión, cantidad"

# Add new row 9: Escuela
$ws.Range("A9").Value = "Escuela"
$ws.Range("B9").Value = "Definición de la escuela"
$ws.Range("C9").Value = "Se calculñaijw0ijfqwopkf"
$ws.Range("D9").Value = "academia"
